# Highlight the remaining "to-do" checklist items in yellow, matching the
# handful of items that were already highlighted before this edit. Each
# target paragraph is matched by its exact text and then the paragraph's
# whole range (text + trailing paragraph mark) is set to wdYellow (7) via
# HighlightColorIndex, which is how Word records <w:highlight w:val="yellow"/>
# on the run(s) inside that paragraph.

$d = $word.ActiveDocument

$targets = @(
    "10 points. Make a logo for your journal blog.",
    "Even though there is minimal need for image editing, you may still need an image editor.",
    "5 points. Find an image that you might use for your blog’s header. Scale and crop this image to be 960 pixels in width and 250 pixels in height. If you get a large photo, you can crop it to this size.",
    "5 points. You must save the image in PNG or JPEG format. It’s best to use JPEG format since most images on Unsplash are photographs.",
    "Store your image in a subfolder named “images”.",
    "5 points. Make the background image of your header tag the image that you found on Unsplash in step 1. This step will require CSS.",
    "Using the “height” property, set the height of your header tag to be the height of your image.",
    "3 points. Change the color of your H1 tag to something that contrasts nicely with your image. Use your judgment in picking the right color (but you must choose something other than the body’s background color).",
    "4 points. Make the header have rounded corners, set the border size to 5 pixels, set the border to the same color as your text (which must be different from your background), set the margins to “auto”, give the text some spacing, and set the width of the header to be something appropriate based on the length of the text.",
    "5 points. Set the width of the main element to 65% of your screen.",
    "5 points. Set the float property to “left”.",
    "3 points. Add an “h2” tag describing the best part of your city.",
    "3 points. Add some dummy text to your main element.",
    "10 points. Add two images. Resize your photos to 800 or fewer pixels wide without distortion. No resizing is necessary if they are already less than 800 pixels wide. Use CSS to scale these images to be 100% width with a 3-pixel padding and a 1-pixel solid border.",
    "3 points. Make this box look nice. I was hoping you could play with CSS and make it look nice. In other words, if you have the bare minimum CSS listed above, you will lose 3 points.",
    "5 points. Set the aside element’s width to 30% of your screen.",
    "5 points. Set the float property to be “right”.",
    "6 points. Add an unordered list of links to your sidebar. Give it the “h2” tag, “More interesting things to see”.",
    "All of these parts will go in the “footer” tag.",
    "5 points. Add “Copyright”, the copyright symbol, the current year, and your name.",
    "5 points. Make this box rest below all other boxes using the “clear” property."
)

$remaining = @{}
foreach ($t in $targets) { $remaining[$t] = $true }

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text
    # Paragraph.Range.Text includes the trailing paragraph mark (CR); strip
    # it off for comparison against the plain checklist text.
    $trimmed = $text.TrimEnd([char]13)
    if ($remaining.ContainsKey($trimmed)) {
        $p.Range.Duplicate.HighlightColorIndex = 7
        $remaining.Remove($trimmed)
    }
}

if ($remaining.Count -gt 0) {
    Write-Host "WARNING: did not find $($remaining.Count) target paragraph(s):"
    foreach ($r in $remaining.Keys) { Write-Host " - $r" }
} else {
    Write-Host "All target paragraphs highlighted."
}
